# Rename the Pearson/PearsonLogo images (footers) from image2.png -> image1.png,
# and the BTec_Logo-Orange image (header) from image1.jpg -> image2.jpg.
#
# These pictures are inline (wp:inline) in the first-page header and in both
# the primary and first-page footers, so walk every section's Headers /
# Footers collections and rename by matching each picture's AlternativeText
# (the descr, which this edit does not change) to decide the new name.

$d = $word.ActiveDocument

function Rename-LogoInRange($rng) {
    if ($rng -eq $null) { return }
    $shapes = $rng.InlineShapes
    if ($shapes -eq $null) { return }
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $descr = $shp.AlternativeText
        if ($descr -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
            $shp.Name = "image1.png"
        } elseif ($descr -eq "BTec_Logo-Orange") {
            $shp.Name = "image2.jpg"
        }
    }
}

for ($s = 1; $s -le $d.Sections.Count; $s++) {
    $sec = $d.Sections.Item($s)

    for ($h = 1; $h -le 3; $h++) {
        $hdr = $sec.Headers.Item($h)
        if ($hdr.Exists) {
            Rename-LogoInRange $hdr.Range
        }
    }

    for ($f = 1; $f -le 3; $f++) {
        $ftr = $sec.Footers.Item($f)
        if ($ftr.Exists) {
            Rename-LogoInRange $ftr.Range
        }
    }
}
